# RF001 - Autenticar Usuario: bump test-suite content from v1.4 to v1.5.
#
# TC1 step 2 (D11) now reports the CAS outage (previously the
# incorrect-username/password alert, which moves to TC2).
# TC2 step 2 (D21) now reports the incorrect-username/password alert
# (previously the TJSeg outage, which moves to TC3), and TC2 step 3 (B22)
# switches to the "preenche os campos" action text.
# TC3 step 2 (D31) now reports the TJSeg outage (previously the CAS
# outage, which moved to TC1); TC3 step 2 (B31) switches to the
# "seleciona um nome de usuario sugerido" action text while TC3 step 3
# (B32) switches to "preenche os campos".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC1 - Expected Result of step 2
$ws.Range("D11").Value = "SYSTEM alerta que o CAS (sistema de autorizacao login-senha) esta fora do ar"

# TC2 - Expected Result of step 2, and Steps text for step 3
$ws.Range("D21").Value = "SYSTEM alerta que o nome de usuario e/ou senha estao incorretos"
$ws.Range("B22").Value = "Usuario do Sistema preenche os campos e clica no botao entrar"

# TC3 - Steps text for step 2, Expected Result of step 2, and Steps text for step 3
$ws.Range("B31").Value = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"
$ws.Range("D31").Value = "SYSTEM alerta que o TJSeg (sistema que fornece as permissoes de acesso e escrita) esta fora do ar"
$ws.Range("B32").Value = "Usuario do Sistema preenche os campos e clica no botao entrar"
